$wb = $excel.ActiveWorkbook

# Locate the Spain worksheet (last sheet) and duplicate it to create the new
# Croatia sheet, placed immediately after Spain (i.e. at the end).
$spain = $wb.Worksheets.Item("Spain")
$spain.Copy([System.Reflection.Missing]::Value, $spain)
$croatia = $wb.Worksheets.Item($wb.Worksheets.Count)
$croatia.Name = "Croatia"

# Fill in the Croatia-specific data. The NGC reference (B4) is entered before
# the market name (B2) so that the new shared-string entries land in the same
# order as the reference workbook (NGC string first, then "Croatia Market").
$croatia.Range("B4").Value = "NGC-3103/T2485/T2488/T2486"
$croatia.Range("B2").Value = "Croatia Market"

# The rows that used to wrap to a taller height on the Spain sheet don't need
# to be as tall for the shorter Croatia text, so re-fit them back down to the
# sheet's normal row height.
$croatia.Range("A3:D5").Rows.AutoFit()

# Column B was sized to fit Spain's longer strings; resize it down for the
# shorter Croatia content.
$croatia.Columns.Item(2).ColumnWidth = 36.5

# Deselect Spain (it was the previously active/selected tab) and select the
# whole sheet range, matching the "not the active tab" selection state used
# by the other, non-active sheets in this workbook.
$spain.Activate() | Out-Null
$spain.Range("A1:XFD1048576").Select() | Out-Null

# Make the new Croatia sheet the active tab, with B9 selected, mirroring the
# usual authoring pattern of landing on the data-entry row after filling in a
# new sheet.
$croatia.Activate() | Out-Null
$croatia.Range("B9").Select() | Out-Null
